# Rename the three logo pictures living in the document's headers/footers.
#
#   header1.xml (First-page header, BTec logo)            image1.jpg -> image2.jpg
#   footer1.xml (First-page footer, Pearson Edexcel logo)  image2.png -> image1.png
#   footer2.xml (Default footer,    Pearson Edexcel logo)  image2.png -> image1.png
#
# These are inline pictures, so they live on the first-page Header
# (wdHeaderFooterFirstPage = 2) and on both the primary Footer
# (wdHeaderFooterPrimary = 1) and the first-page Footer
# (wdHeaderFooterFirstPage = 2) of the document's one section.

$d = $word.ActiveDocument
$sec = $d.Sections.First

# First-page header: BTec_Logo-Orange -> rename image1.jpg to image2.jpg
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists) {
    for ($i = 1; $i -le $hdrFirst.Range.InlineShapes.Count; $i++) {
        $shp = $hdrFirst.Range.InlineShapes.Item($i)
        $shp.Name = "image2.jpg"
    }
}

# Primary footer: PearsonLogo -> rename image2.png to image1.png
$ftPrimary = $sec.Footers.Item(1)
if ($ftPrimary.Exists) {
    for ($i = 1; $i -le $ftPrimary.Range.InlineShapes.Count; $i++) {
        $shp = $ftPrimary.Range.InlineShapes.Item($i)
        $shp.Name = "image1.png"
    }
}

# First-page footer: PearsonLogo -> rename image2.png to image1.png
$ftFirst = $sec.Footers.Item(2)
if ($ftFirst.Exists) {
    for ($i = 1; $i -le $ftFirst.Range.InlineShapes.Count; $i++) {
        $shp = $ftFirst.Range.InlineShapes.Item($i)
        $shp.Name = "image1.png"
    }
}

Write-Output "Renamed header/footer logo pictures."
